$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two "type_respiro" category labels used in column H
# ("A" -> "HI" and "B" -> "LI") across all data rows.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 17 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $v = $cell.Value2
    if ($v -eq "A") {
        $cell.Value = "HI"
    } elseif ($v -eq "B") {
        $cell.Value = "LI"
    }
}

# Update the current selection to match the saved view state
[void]$ws.Range("F23").Select()
